# Rotate the data of rows 26, 27, 28 on the active sheet:
#   new row26 <- old row28
#   new row27 <- old row26
#   new row28 <- old row27
# (Columns K, L, M, N, AC are handled separately below because they are
# only populated on the row that used to be row 26.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","P","Q","R","S","T","U","V","W","Y","AA","AD","AE","AG","AW","AX")

# Columns whose text looks like a date ("2023-09-14") must be forced to stay
# text instead of being auto-parsed into a date serial number when written
# back through .Value2.
$textCols = @("Y","AA")

# Snapshot the current (pre-edit) values for rows 26-28, per column.
$row26 = @{}
$row27 = @{}
$row28 = @{}
foreach ($col in $cols) {
    $row26[$col] = $ws.Range("$col" + "26").Value2
    $row27[$col] = $ws.Range("$col" + "27").Value2
    $row28[$col] = $ws.Range("$col" + "28").Value2
}

# Apply the rotation: 26<-28, 27<-26, 28<-27
foreach ($col in $cols) {
    $isText = $textCols -contains $col
    if ($isText) { $ws.Range("$col" + "26").NumberFormat = "@" }
    $ws.Range("$col" + "26").Value2 = $row28[$col]
    if ($isText) { $ws.Range("$col" + "26").Style = "Normal" }

    if ($isText) { $ws.Range("$col" + "27").NumberFormat = "@" }
    $ws.Range("$col" + "27").Value2 = $row26[$col]
    if ($isText) { $ws.Range("$col" + "27").Style = "Normal" }

    if ($isText) { $ws.Range("$col" + "28").NumberFormat = "@" }
    $ws.Range("$col" + "28").Value2 = $row27[$col]
    if ($isText) { $ws.Range("$col" + "28").Style = "Normal" }
}

# Columns K, L, M, N and AC only have content on one of these three rows
# (the row that used to be row 26). Move that marker along with its row,
# i.e. onto the row that is now 27.
$ws.Range("K26").ClearContents()
$ws.Range("L26").ClearContents()
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("AC26").ClearContents()

# Plain "" writes are dropped entirely by the engine (no value == no cell),
# so round-trip these through a temporary text number-format to force the
# cell to persist as an (empty) entry, then restore the default style.
foreach ($col in @("K","L","M","N")) {
    $r = $ws.Range("$col" + "27")
    $r.NumberFormat = "@"
    $r.Value2 = ""
    $r.Style = "Normal"
}
$ws.Range("AC27").Value2 = "ringhack, äldre"

$ws.Range("K28").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("AC28").ClearContents()
